$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "92.131.98"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "3.101.72"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "614.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.40%  "
$ws.Range("E7").Value = "  -5.85%  "
$ws.Range("E8").Value = "  +4.89%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("D10").Value = "3.099.67"
$ws.Range("E10").Value = "  -0.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.729"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.67%  "
$ws.Range("E12").Value = "  -1.65%  "
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.37%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "91.735.84"
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").Value = "3.668.99"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").Value = "3.087.65"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "447.60"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.78%  "
$ws.Range("E24").Value = "  -5.08%  "
$ws.Range("E25").Value = "  -1.62%  "
$ws.Range("B26").Value = "Aptos"
$ws.Range("C26").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.48%  "
$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "81.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.70%  "
$ws.Range("D28").Value = "3.264.69"
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.137"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +13.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.236"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.167"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +57.15%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.94"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.98%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.163"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.84%  "
$ws.Range("E37").Value = "  -2.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.14"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.92"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "479.98"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.36%  "
$ws.Range("E41").Value = "  -0.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.432"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.75%  "
$ws.Range("E43").Value = "  -4.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "159.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.36%  "
$ws.Range("E48").Value = "  -1.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0332"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.79%  "
